$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Problem 2.2 - Base case")

# Generator Data table
$ws.Range("C4").Value = 50
$ws.Range("K4").Value = 50
$ws.Range("Q4").Value = 100
$ws.Range("R4").Value = -19.84

$ws.Range("C5").Value = 45
$ws.Range("D5").Value = "Node 4"
$ws.Range("K5").Value = 170
$ws.Range("Q5").Value = 10000
$ws.Range("R5").Value = -26.88

# Remove Gen 3 (row 6, columns A:E)
$ws.Range("A6:E6").ClearContents()

$ws.Range("K6").Value = 200
$ws.Range("P6").Value = "Line 2-4"
$ws.Range("Q6").Value = 10000
$ws.Range("R6").Value = -26.88

# New row 7
$ws.Range("K7").Value = 80
$ws.Range("L7").Value = "Node 4"
$ws.Range("P7").Value = "Line 3-4"
$ws.Range("Q7").Value = 100
$ws.Range("R7").Value = -15.72
$ws.Range("J7").Value = "Load 4"

$ws.Range("G23").Select()
